# Update "想去人数" (F column) counts on sheet "展览" and sheet "全部类型"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 63
$ws1.Range("F6").Value = 181
$ws1.Range("F7").Value = 3766
$ws1.Range("F8").Value = 182
$ws1.Range("F11").Value = 75
$ws1.Range("F12").Value = 80
$ws1.Range("F13").Value = 664
$ws1.Range("F14").Value = 161
$ws1.Range("F15").Value = 893
$ws1.Range("F16").Value = 67
$ws1.Range("F20").Value = 89
$ws1.Range("F22").Value = 3248
$ws1.Range("F23").Value = 5592
$ws1.Range("F27").Value = 508
$ws1.Range("F29").Value = 3200
$ws1.Range("F30").Value = 341
$ws1.Range("F31").Value = 2399
$ws1.Range("F35").Value = 176
$ws1.Range("F37").Value = 339
$ws1.Range("F39").Value = 493
$ws1.Range("F40").Value = 870
$ws1.Range("F42").Value = 26
$ws1.Range("F44").Value = 55
$ws1.Range("F45").Value = 532

$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F4").Value = 63
$ws2.Range("F6").Value = 181
$ws2.Range("F7").Value = 3766
$ws2.Range("F8").Value = 182
$ws2.Range("F12").Value = 75
$ws2.Range("F13").Value = 80
$ws2.Range("F14").Value = 664
$ws2.Range("F15").Value = 161
$ws2.Range("F16").Value = 893
$ws2.Range("F17").Value = 67
$ws2.Range("F21").Value = 89
$ws2.Range("F23").Value = 3249
$ws2.Range("F24").Value = 5592
$ws2.Range("F28").Value = 508
$ws2.Range("F30").Value = 3200
$ws2.Range("F31").Value = 341
$ws2.Range("F32").Value = 2399
$ws2.Range("F35").Value = 111
$ws2.Range("F36").Value = 176
$ws2.Range("F38").Value = 339
$ws2.Range("F40").Value = 493
$ws2.Range("F41").Value = 870
$ws2.Range("F43").Value = 26
$ws2.Range("F45").Value = 55
$ws2.Range("F46").Value = 532
